$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:T9").ClearContents()
$ws.Rows("6:9").Delete()

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Inha"
$ws.Range("C2").Value = "Tgfbr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.390455
$ws.Range("H2").Value = 1.171365
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 42.09975866666667
$ws.Range("N2").Value = 126.299276
$ws.Range("O2").Value = 0.3315552933456474
$ws.Range("P2").Value = 0.3315552933456474
$ws.Range("Q2").Value = 16.43806127019333
$ws.Range("R2").Value = 147.94255143174
$ws.Range("S2").Value = 0.3315552933456474
$ws.Range("T2").Value = 0.3315552933456474

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Inha"
$ws.Range("C3").Value = "Tgfbr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.390455
$ws.Range("H3").Value = 1.171365
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 57.16769933333334
$ws.Range("N3").Value = 171.503098
$ws.Range("O3").Value = 0.4502223747274475
$ws.Range("P3").Value = 0.4502223747274475
$ws.Range("Q3").Value = 22.32141404319667
$ws.Range("R3").Value = 200.89272638877
$ws.Range("S3").Value = 0.4502223747274475
$ws.Range("T3").Value = 0.4502223747274475

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Inha"
$ws.Range("C4").Value = "Tgfbr3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.390455
$ws.Range("H4").Value = 1.171365
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 27.596267
$ws.Range("N4").Value = 82.78880100000001
$ws.Range("O4").Value = 0.2173335118824389
$ws.Range("P4").Value = 0.2173335118824389
$ws.Range("Q4").Value = 10.775100431485
$ws.Range("R4").Value = 96.97590388336501
$ws.Range("S4").Value = 0.2173335118824389
$ws.Range("T4").Value = 0.2173335118824389

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Inha"
$ws.Range("C5").Value = "Tgfbr3"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.390455
$ws.Range("H5").Value = 1.171365
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1128593333333333
$ws.Range("N5").Value = 0.338578
$ws.Range("O5").Value = 0.0008888200444663087
$ws.Range("P5").Value = 0.0008888200444663087
$ws.Range("Q5").Value = 0.04406649099666667
$ws.Range("R5").Value = 0.39659841897
$ws.Range("S5").Value = 0.0008888200444663087
$ws.Range("T5").Value = 0.0008888200444663087
